$d = $word.ActiveDocument

# Remove the 6 leading empty paragraphs.
for ($i = 0; $i -lt 6; $i++) {
    $d.Paragraphs(1).Range.Delete()
}
